# Apply edits described by the diff:
# 1. Update shared strings "VP-NCC-R-004-00X" -> "VP-NCC-S-013-00X" (cells A1:A5)
# 2. Update the active cell selection on the sheet from C11 to B11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "VP-NCC-S-013-001"
$ws.Range("A2").Value = "VP-NCC-S-013-002"
$ws.Range("A3").Value = "VP-NCC-S-013-003"
$ws.Range("A4").Value = "VP-NCC-S-013-004"
$ws.Range("A5").Value = "VP-NCC-S-013-005"

$ws.Activate()
$ws.Range("B11").Select()
